# Updated cryptos list on Sun Jul  2 19:32:18 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.537.37"
$ws.Range("E2").Value = "  -0.30%  "

# Row 3
$ws.Range("D3").Value = "1.911.39"
$ws.Range("E3").Value = "  -0.62%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.24"
$ws.Range("E5").Value = "  -1.19%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("E6").Value = "  -0.07%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4839"
$ws.Range("E7").Value = "  +1.99%  "

# Row 8
$ws.Range("E8").Value = "  +0.11%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06817"
$ws.Range("E9").Value = "  -0.35%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "111.27"
$ws.Range("E10").Value = "  +5.78%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.26"
$ws.Range("E11").Value = "  +4.76%  "

# Row 12
$ws.Range("D12").Value = "1.916.25"
$ws.Range("E12").Value = "  -0.29%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07566"
$ws.Range("E13").Value = "  -1.71%  "

# Row 14
$ws.Range("E14").Value = "  +1.09%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6703"
$ws.Range("E15").Value = "  +0.34%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "297.19"
$ws.Range("E16").Value = "  +1.88%  "

# Row 17
$ws.Range("D17").Value = "30.528.26"
$ws.Range("E17").Value = "  -0.34%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.02"
$ws.Range("E18").Value = "  +0.48%  "

# Row 19 (name/link/price/volume swap)
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007600"
$ws.Range("E19").Value = "  -0.31%  "

# Row 20 (name/link/price/volume swap)
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9998"
$ws.Range("E20").Value = "  +0.01%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.521"
$ws.Range("E21").Value = "  -1.24%  "

# Row 22
$ws.Range("D22").Value = "2.165.87"
$ws.Range("E22").Value = "  -0.26%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9993"
$ws.Range("E23").Value = "  -0.08%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.428"
$ws.Range("E24").Value = "  -0.14%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.465"
$ws.Range("E25").Value = "  +0.11%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.72"
$ws.Range("E26").Value = "  -1.27%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.34"
$ws.Range("E27").Value = "  -3.65%  "

# Row 28
$ws.Range("E28").Value = "  -1.75%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1065"
$ws.Range("E29").Value = "  -0.76%  "

# Row 30
$ws.Range("E30").Value = "  +2.91%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.143"
$ws.Range("E31").Value = "  -0.96%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.047"
$ws.Range("E32").Value = "  -0.45%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04981"
$ws.Range("E33").Value = "  -1.23%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7367"
$ws.Range("E34").Value = "  -0.24%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.136"
$ws.Range("E35").Value = "  -0.82%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9994"
$ws.Range("E36").Value = "  +0.04%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02040"
$ws.Range("E37").Value = "  -1.55%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.714"
$ws.Range("E38").Value = "  -0.86%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.679"
$ws.Range("E39").Value = "  -0.46%  "

# Row 40
$ws.Range("E40").Value = "  -1.95%  "

# Row 41 (name/link/price/volume swap)
$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "109.31"
$ws.Range("E41").Value = "  -1.72%  "

# Row 42 (name/link/price/volume swap)
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4452"
$ws.Range("E42").Value = "  +1.41%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8676"
$ws.Range("E43").Value = "  -1.03%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.787"
$ws.Range("E44").Value = "  -1.94%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9992"
$ws.Range("E45").Value = "  -0.08%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "69.45"
$ws.Range("E46").Value = "  +2.00%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.213"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.30"
$ws.Range("E48").Value = "  -0.02%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.240"
$ws.Range("E49").Value = "  -0.79%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1227"
$ws.Range("E50").Value = "  -1.31%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2513"
$ws.Range("E51").Value = "  -0.37%  "
